# Refresh the Nanning comic-expo listing to the latest scrape.
#
# The oldest event that already took place
# ("南宁·原神x星铁x绝区零同人ONLY3.0", 2024-09-15) has aged out of the feed,
# so its row is removed from both the "展览" (exhibitions) and "全部类型"
# (all types) sheets. Every later row shifts up one slot, the running
# index kept in column A is renumbered to stay sequential, and two
# "想去人数" (want-to-go) counters that were simply refreshed between
# scrapes (for events whose other details did not change) get their
# updated values.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the last used row before we touch anything.
    $lastRowBefore = $ws.UsedRange.Rows.Count

    # Drop the row for the expo that has already ended (old row 2); this
    # shifts every following row up by one, the same as the source feed
    # dropping the stale entry.
    $ws.Rows.Item(2).Delete()

    $lastRowAfter = $lastRowBefore - 1

    # Row deletion leaves column A's stored numbers untouched (merely
    # shifted up), so renumber the running index (1, 2, 3, ...) for every
    # remaining data row.
    for ($r = 2; $r -le $lastRowAfter; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # "想去人数" (want-to-go count) was refreshed between scrapes for two
    # events whose other details otherwise just shifted up unchanged.
    $ws.Range("F3").Value = 298
    $ws.Range("F4").Value = 4319
}
